$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 2-6 down to 3-7 to open a slot for the new vacancy at row 2,
# then re-assert every cell so the final row order matches the production data
# (row order became 2(new),4,5,3,6 relative to the original rows after the new entry was added).
$ws.Rows("2:2").Insert()

# ---- Row 2 ----
$ws.Range("A2").Value = "Acompañante de menor o Maestro Sombra"
$ws.Range("B2").Value = "La empresa es confidencial o no se encuentra disponible"
$ws.Range("C2").Value = "Hermosillo,, Son."
$ws.Range("D2").Value = "`$12,000 Mensual"
$ws.Range("E2").Value = $True
$ws.Range("F2").Value = "Educación"
$ws.Range("G2").Value = "Educación especial"
$ws.Range("H2").Value = "Universitario titulado"
$ws.Range("I2").Value = "Permanente"
$ws.Range("J2").Value = "Tiempo completo"
$ws.Range("K2").Value = "Presencial"
$descripcion2 = @'
Vacante: Acompañante de niño (TDA o Autismo) – Maestro Sombra
Ubicación: Zona Poniente, Hermosillo
Horario: lunes a viernes
Requisitos:
Sexo indistinto
Carrera en Pedagogía, Ciencias de la Educación o afín
Especialidad o experiencia en Educación Especial
Experiencia como maestro sombra o acompañante educativo
Paciencia, empatía y habilidades de comunicación
Responsabilidad y compromiso
Preferente cuente con carro para su traslado
Funciones principales:
Brindar apoyo personalizado a un niño con TDA o Trastorno del Espectro Autista en su entorno escolar
Favorecer la integración e inclusión en actividades académicas y sociales
Implementar estrategias de apoyo de acuerdo con el plan educativo
Colaborar con docentes y padres para dar seguimiento al progreso
Ofrecemos:
Contratación directa
Estabilidad laboral
Ambiente de trabajo respetuoso y colaborativo
'@
$ws.Range("L2").Value = $descripcion2

# ---- Row 3 ----
$ws.Range("A3").Value = "Terapeuta"
$ws.Range("B3").Value = "La empresa es confidencial o no se encuentra disponible"
$ws.Range("C3").Value = "León,, Gto."
$ws.Range("D3").Value = "`$14,000 - `$16,000 Mensual"
$ws.Range("E3").Value = $True
$ws.Range("F3").Value = "Sector salud"
$ws.Range("G3").Value = "Terapeuta"
$ws.Range("H3").Value = "Universitario titulado"
$ws.Range("I3").Value = "Permanente"
$ws.Range("J3").Value = "Tiempo completo"
$ws.Range("K3").Value = "Presencial"
$descripcion3 = @'
Requisitos del puesto




Estudios universitarios con título en Terapia.
Experiencia previa como Terapeuta de niños con Trastornos del espectro autista.
Gusto por realizar manualidades.
Habilidad para nadar.
Licencia de manejo vigente.
Conocimientos en técnicas de terapia y rehabilitación.
Licencia o certificación válida en Terapia (deseable).







Responsabilidades del puesto




Realizar evaluaciones y diagnósticos de los pacientes.
Diseñar planes de tratamiento personalizados.
Realizar sesiones de terapia adaptadas a las necesidades individuales de cada paciente.
Mantener registros precisos de la evolución de los pacientes.







Prestaciones y beneficios adicionales




Salario mensual competitivo de 14000 a 16000.
Prestaciones de ley.
Vales de despensa.
Fondo de ahorro.
Oportunidades de capacitación y desarrollo profesional en un ambiente de trabajo colaborativo y en constante crecimiento.
'@
$ws.Range("L3").Value = $descripcion3

# ---- Row 4 ----
$ws.Range("A4").Value = "Psicología clínica"
$ws.Range("B4").Value = "Fundación Planeta Tea AC"
$ws.Range("C4").Value = "San Pedro Garza García,, N.L."
$ws.Range("D4").Value = "`$14,000 - `$15,000 Mensual"
$ws.Range("E4").Value = $True
$ws.Range("F4").Value = "Ciencias sociales - Humanidades"
$ws.Range("G4").Value = "Psicología"
$ws.Range("H4").Value = "Universitario titulado"
$ws.Range("I4").Value = "Permanente"
$ws.Range("J4").Value = "Tiempo completo"
$ws.Range("K4").Value = "Presencial"
$descripcion4 = @'
Acerca de la empresa




Fundación Planeta Tea AC es una organización dedicada a apoyar a niños y adultos con autismo de escasos recursos y sus familias, brindando servicios de psicología clínica y terapias especializadas. Nuestra misión es mejorar la calidad de vida de las personas con autismo a través de la educación y la inclusión en la sociedad. Al unirte a nuestro equipo, tendrás la oportunidad de contribuir a una causa noble y en constante crecimiento, además de recibir apoyo y capacitación continua. Ubicación: San Pedro Garza García, Nuevo León.




Requisitos del puesto







Educación mínima: Universitario con título y cédula profesional en Psicología
Experiencia y conocimiento en la aplicación de pruebas psicológicas.
Experiencia en pacientes con autismo de mínimo 1 año y 3 años deseable.







Responsabilidades del puesto




Evaluación emocional y conductual
Psicoterapia individual y familiar
Apoyo en planes educativos individualizados.
Manejo de citas
Coordinación de personal







Prestaciones y beneficios adicionales




Prestaciones de ley
Salario mensual de 14000 a 15000.
Bonos incentivos
Horario a tiempo completo: Lunes a viernes de 10am a 7pm, con 1 hora de comida y Sábados de 9am a 1:30pm horario corrido.
Modalidad presencial.
Oportunidades de capacitación y desarrollo profesional.
Excelente ambiente de trabajo.
'@
$ws.Range("L4").Value = $descripcion4

# ---- Row 5 ----
$ws.Range("A5").Value = "Monitora en inclusión educativa (maestro sombra)"
$ws.Range("B5").Value = "CC INTEGRACION LABORAL"
$ws.Range("C5").Value = "CDMX"
$ws.Range("D5").Value = "`$8,500 Mensual"
$ws.Range("E5").Value = $False
$ws.Range("F5").Value = "Educación"
$ws.Range("G5").Value = "Educación especial"
$ws.Range("H5").Value = "Universitario sin titulo"
$ws.Range("I5").Value = "Permanente"
$ws.Range("J5").Value = "Tiempo completo"
$ws.Range("K5").Value = "Presencial"
$descripcion5 = @'
REQUISITOS:

Nivel de estudios: Licenciatura (concluida o últimos semestres) Psicología educativa, Pedagogía, Educación especial o afines

23 a 32 años

Sexo indistinto

Estado civil indistinto




EXPERIENCIA:

Deseable con niños con autismo, trastornos del neurodesarrollo y/o alguna discapacidad




HABILIDADES Y COMPETENCIAS:

Proactiva, comunicación asertiva, responsable, puntual.

Empatía y sensibilidad emocional, paciencia, tolerancia a la frustración, vocación infantil, compromiso y responsabilidad.




HORARIO DE TRABAJO:

De lunes a viernes




Interesados enviar cv a la dirección de contacto.
'@
$ws.Range("L5").Value = $descripcion5

# ---- Row 6 ----
$ws.Range("A6").Value = "Psicóloga"
$ws.Range("B6").Value = "Neuro Activa"
$ws.Range("C6").Value = "CDMX"
$ws.Range("D6").Value = "`$8,364 - `$8,500 Mensual"
$ws.Range("E6").Value = $False
$ws.Range("F6").Value = "Educación"
$ws.Range("G6").Value = "Psicopedagogía"
$ws.Range("H6").Value = "Universitario titulado"
$ws.Range("I6").Value = "Permanente"
$ws.Range("J6").Value = "Tiempo completo"
$ws.Range("K6").Value = "Híbrido"
$descripcion6 = @'
Acerca de la empresa




Neuro Activa es una empresa líder en el campo de la educación y la psicopedagogía, comprometida con el desarrollo integral de cada individuo. Trabajar con nosotros significa formar parte de un equipo dedicado a brindar soluciones innovadoras y personalizadas a las necesidades educativas especiales. - Ubicación: Ciudad de México.




Requisitos del puesto




Licenciatura en Psicología, Educación Especial, Pedagogía, Psicopedagogía o afines.
Diplomado o curso en intervención educativa, inclusión escolar o TEA/TDAH. (Deseable)
Experiencia mínima de 6 meses con niños con necesidades educativas especiales (NEE), Trastorno del Espectro Autista (TEA), TDAH u otras condiciones del neurodesarrollo.
Conocimientos básicos de adaptaciones curriculares y manejo de conductas.
Habilidades en regulación emocional, contención y acompañamiento respetuoso.
Capacidad para implementar estrategias sensoriales, estructuración del entorno y apoyos visuales.
Buen manejo de la comunicación con familia y equipo docente.







Responsabilidades del puesto




Brindar intervención educativa y apoyo a niños con necesidades educativas especiales.
Implementar estrategias de inclusión escolar y adaptaciones curriculares.
Colaborar con el equipo docente y familia en el desarrollo del plan de intervención individualizado.













Únete a nuestro equipo en Neuro Activa y sé parte de un proyecto que transforma vidas a través de la educación inclusiva. ¡Aplica ahora y haz la diferencia!
'@
$ws.Range("L6").Value = $descripcion6

# ---- Row 7 ----
$ws.Range("A7").Value = "Maestra"
$ws.Range("B7").Value = "Asociación Centro de T..."
$ws.Range("C7").Value = "SLP."
$ws.Range("D7").Value = "`$5,000 - `$6,000 Mensual"
$ws.Range("E7").Value = $False
$ws.Range("F7").Value = "Educación"
$ws.Range("G7").Value = "Educación especial"
$ws.Range("H7").Value = "Universitario titulado"
$ws.Range("I7").Value = "Permanente"
$ws.Range("J7").Value = "Medio tiempo"
$ws.Range("K7").Value = "Híbrido"
$descripcion7 = @'
Acerca de la empresa




**** Asociación Centro de Terapia Infantil y de Educación Especial es una organización comprometida con el bienestar y desarrollo de niños, niñas y jóvenes con autismo, síndrome de Down y/o discapacidad intelectual. Nuestra misión es brindar un ambiente inclusivo y terapéutico para promover su crecimiento y aprendizaje. - **Ubicación:** San Luis Potosí.




Requisitos del puesto




****
Educación mínima requerida: Universitario con título.
Habilidades: Paciente y tolerante, proactiva y dinámica, ordenada y puntual.
Experiencia: Experiencia previa en educación especial es deseable.







Responsabilidades del puesto




****
Diseñar e implementar actividades educativas y terapéuticas personalizadas.
Colaborar con el equipo multidisciplinario para el seguimiento y evaluación del progreso de los estudiantes.
Mantener un ambiente seguro y estimulante para los niños y jóvenes atendidos.







Prestaciones y beneficios adicionales




****
Capacitaciones continuas en educación especial.
Oportunidades de crecimiento profesional.
Ambiente de trabajo colaborativo y respetuoso.
'@
$ws.Range("L7").Value = $descripcion7
